# Update diary as of 13 Feb 2020
# Fills in the two diary entries for 12 Jan 2020 (W) and 13 Jan 2020 (Th)
# in rows 33 and 34 of Sheet1, which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: 12 Jan 2020 (W)
$ws.Range("A33").Value = "12 Jan 2020 (W)"
$ws.Range("B33").Value = "2245-0000"
$ws.Range("C33").Value = "Harry"
$ws.Range("D33").Value = "Study for tomorrow’s midterm"
$ws.Range("E33").Value = "Got through the lectures and online articles"
$ws.Range("F33").Value = "I don’t procrastinate out of laziness; the past week or so were full of chaos, both in and out of academia. On the bright side, the lecture recordings help tremendously with the studying."
$ws.Range("G33").Value = "Feeling fine, just need a smoothie. I’ll get one on the way to class tomorrow."

# Row 34: 13 Jan 2020 (Th)
$ws.Range("A34").Value = "13 Jan 2020 (Th)"
$ws.Range("B34").Value = "0000-0206"
$ws.Range("C34").Value = "Harry"
$ws.Range("D34").Value = "Study for today’s midterm"
$ws.Range("E34").Value = "Same as yesterday. Will finish studying later today."
$ws.Range("F34").Formula = "=F33"
$ws.Range("G34").Value = "Stomach is active, just need a smoothie. I’ll get one on the way to class today."

# Both rows grow to fit the wrapped reflection text (same height as similar
# multi-line entries elsewhere in the diary).
$ws.Rows.Item(33).RowHeight = 73.1
$ws.Rows.Item(34).RowHeight = 73.1

# Move the selection to G34, mirroring where editing left off.
$ws.Range("G34").Select()
